$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: column D numeric-looking strings must stay as text
# (match the source price-string formatting, e.g. trailing zeros / leading zeros)
# are set via NumberFormat "@" before assignment so Excel does not coerce them to numbers.
$textCells = @(
    "D5", "D6", "D7", "D8", "D10", "D11", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values (coin rankings refresh).
$ws.Range("D2").Value = "67.274.49"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "3.342.19"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "578.39"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "184.15"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.604"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "0.406"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "3.926.21"
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "27.34"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "67.490.27"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "3.346.56"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "445.35"
$ws.Range("E18").Value = "  +6.94%  "
$ws.Range("D19").Value = "13.61"
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").Value = "5.66"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "7.72"
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("D22").Value = "74.02"
$ws.Range("E22").Value = "  +4.00%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "3.490.70"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").Value = "0.512"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "0.0000120"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").Value = "0.195"
$ws.Range("E27").Value = "  +4.41%  "
$ws.Range("D28").Value = "9.06"
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "1.97"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").Value = "22.93"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").Value = "5.33"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "6.78"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +4.47%  "
$ws.Range("D37").Value = "161.97"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").Value = "27.53"
$ws.Range("E38").Value = "  +3.78%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").Value = "2.835.15"
$ws.Range("E40").Value = "  +8.13%  "
$ws.Range("D41").Value = "0.792"
$ws.Range("D42").Value = "4.47"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "40.32"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "0.0672"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "24.52"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "2.36"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("D48").Value = "323.25"
$ws.Range("E48").Value = "  -3.42%  "
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").Value = "0.987"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").Value = "31.00"
$ws.Range("E51").Value = "  +1.99%  "
